# daily auto push: 2026-02-19 19:05 UTC
#
# Inserts one new data row at row 843 (2026/02/19, 木, 23, 201) into
# Sheet1's log table, pushing the former rows 843..884 down to 844..885.
# This expands the used range from A1:D884 to A1:D885.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 843..884 down by one, leaving a blank row 843 to fill in.
$ws.Rows.Item(843).Insert()

# Column A holds dates formatted/stored as plain text (e.g. "2026/02/19"),
# not real date serials. Force text formatting before assigning the value
# so Excel doesn't auto-convert the string into a date serial number, then
# drop the formatting again so the cell ends up with the sheet's default
# (unstyled) look, matching every other data row in the column.
$cellA = $ws.Cells.Item(843, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/02/19"
$cellA.ClearFormats()

$ws.Cells.Item(843, 2).Value = "木"
$ws.Cells.Item(843, 3).Value = 23
$ws.Cells.Item(843, 4).Value = 201
